$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price / 1h-volume figures from the latest symbol-list refresh.
# Cells are stored as text (matching the source data format), so force the
# Text number format before writing to avoid Excel auto-converting the
# numeric-looking / percent-looking strings into Number/Percentage cells.
$cellRefs = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "E26",
    "D27",
    "E27",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "E48",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51",
)
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "303.69"
$ws.Range("E2").Value = "5.36%"
$ws.Range("D3").Value = "34.90"
$ws.Range("E3").Value = "12.21%"
$ws.Range("D4").Value = "5.204"
$ws.Range("E4").Value = "5.55%"
$ws.Range("D5").Value = "0.07799"
$ws.Range("E5").Value = "6.39%"
$ws.Range("D6").Value = "2.373"
$ws.Range("E6").Value = "7.48%"
$ws.Range("D7").Value = "8.042"
$ws.Range("E7").Value = "4.17%"
$ws.Range("D8").Value = "0.9325"
$ws.Range("E8").Value = "3.09%"
$ws.Range("D9").Value = "0.1013"
$ws.Range("E9").Value = "10.55%"
$ws.Range("D10").Value = "0.1846"
$ws.Range("E10").Value = "9.17%"
$ws.Range("D11").Value = "0.08646"
$ws.Range("E11").Value = "5.29%"
$ws.Range("D12").Value = "0.03320"
$ws.Range("E12").Value = "6.57%"
$ws.Range("D13").Value = "0.09902"
$ws.Range("E13").Value = "-0.37%"
$ws.Range("D14").Value = "0.001484"
$ws.Range("E14").Value = "-1.10%"
$ws.Range("D15").Value = "0.005751"
$ws.Range("E15").Value = "0.72%"
$ws.Range("D16").Value = "3.472"
$ws.Range("E16").Value = "-1.48%"
$ws.Range("D17").Value = "3.964"
$ws.Range("E17").Value = "6.23%"
$ws.Range("D18").Value = "2.148"
$ws.Range("E18").Value = "4.99%"
$ws.Range("D19").Value = "0.3373"
$ws.Range("E19").Value = "1.18%"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").Value = "0.34%"
$ws.Range("D21").Value = "4.299"
$ws.Range("E21").Value = "2.88%"
$ws.Range("E22").Value = "6.02%"
$ws.Range("D23").Value = "0.04574"
$ws.Range("E23").Value = "0.89%"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").Value = "0.78%"
$ws.Range("D25").Value = "0.004430"
$ws.Range("E25").Value = "6.46%"
$ws.Range("E26").Value = "-0.11%"
$ws.Range("D27").Value = "0.0003700"
$ws.Range("E27").Value = "8.93%"
$ws.Range("D39").Value = "0.01770"
$ws.Range("E39").Value = "13.21%"
$ws.Range("D40").Value = "0.04809"
$ws.Range("E40").Value = "8.35%"
$ws.Range("D41").Value = "0.007752"
$ws.Range("E41").Value = "5.50%"
$ws.Range("D42").Value = "0.1408"
$ws.Range("E42").Value = "5.91%"
$ws.Range("D43").Value = "0.007134"
$ws.Range("E43").Value = "-25.27%"
$ws.Range("D44").Value = "0.002289"
$ws.Range("E44").Value = "3.04%"
$ws.Range("D45").Value = "0.009199"
$ws.Range("E45").Value = "1.11%"
$ws.Range("D46").Value = "0.00005954"
$ws.Range("E46").Value = "-2.68%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.10%"
$ws.Range("E48").Value = "13.87%"
$ws.Range("E49").Value = "-0.10%"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.10%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.10%"
